$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to text
# format (matching the source file where all these cells are stored as
# inline/shared strings), then the style is reset to Normal so no stray
# number-format style gets attached to the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "30.094.25"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "1.919.30"
$ws.Range("E3").Value = "  +2.46%  "
$ws.Range("E4").Value = "  +0.15%  "
Set-TextValue $ws.Range("D5") "319.26"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("E6").Value = "  +0.14%  "
Set-TextValue $ws.Range("D7") "0.5075"
$ws.Range("E7").Value = "  -0.19%  "
Set-TextValue $ws.Range("D8") "0.4025"
$ws.Range("E8").Value = "  +1.81%  "
Set-TextValue $ws.Range("D9") "0.08318"
$ws.Range("E9").Value = "  +1.51%  "
Set-TextValue $ws.Range("D10") "1.114"
$ws.Range("E10").Value = "  +1.74%  "
Set-TextValue $ws.Range("D11") "42.13"
Set-TextValue $ws.Range("D12") "24.10"
$ws.Range("E12").Value = "  +1.16%  "
Set-TextValue $ws.Range("D13") "6.416"
$ws.Range("E13").Value = "  +1.75%  "
$ws.Range("D14").Value = "1.917.18"
$ws.Range("E14").Value = "  +2.35%  "
Set-TextValue $ws.Range("D15") "7.244"
$ws.Range("E15").Value = "  +0.57%  "
Set-TextValue $ws.Range("D16") "1.001"
$ws.Range("E16").Value = "  +0.01%  "
Set-TextValue $ws.Range("D17") "92.48"
$ws.Range("E18").Value = "  +0.82%  "
Set-TextValue $ws.Range("D19") "0.06492"
$ws.Range("E19").Value = "  +1.56%  "
Set-TextValue $ws.Range("D20") "18.44"
$ws.Range("E20").Value = "  +2.17%  "
$ws.Range("E21").Value = "  +0.17%  "
Set-TextValue $ws.Range("D22") "5.945"
$ws.Range("E22").Value = "  +1.81%  "
$ws.Range("D23").Value = "30.096.70"
$ws.Range("E23").Value = "  +0.27%  "
Set-TextValue $ws.Range("D24") "11.34"
$ws.Range("E24").Value = "  +1.76%  "
Set-TextValue $ws.Range("D25") "2.195"
$ws.Range("E25").Value = "  +0.95%  "
$ws.Range("D26").Value = "2.137.03"
$ws.Range("E26").Value = "  +2.47%  "
Set-TextValue $ws.Range("D27") "21.81"
$ws.Range("E27").Value = "  +2.72%  "
Set-TextValue $ws.Range("D28") "162.47"
$ws.Range("E28").Value = "  +1.08%  "
Set-TextValue $ws.Range("D29") "2.270"
$ws.Range("E29").Value = "  +1.47%  "
Set-TextValue $ws.Range("D30") "129.04"
$ws.Range("E30").Value = "  +1.20%  "
Set-TextValue $ws.Range("D31") "1.133"
$ws.Range("E31").Value = "  +5.71%  "
Set-TextValue $ws.Range("D32") "0.1046"
$ws.Range("E32").Value = "  +1.04%  "
Set-TextValue $ws.Range("D33") "5.939"
$ws.Range("E33").Value = "  -0.43%  "
Set-TextValue $ws.Range("D34") "3.783"
$ws.Range("E34").Value = "  +1.77%  "
Set-TextValue $ws.Range("D35") "0.02449"
$ws.Range("E35").Value = "  +0.30%  "
Set-TextValue $ws.Range("D36") "5.305"
$ws.Range("E36").Value = "  +1.24%  "
Set-TextValue $ws.Range("D37") "0.06443"
$ws.Range("E37").Value = "  +0.77%  "
Set-TextValue $ws.Range("D38") "1.236"
$ws.Range("E38").Value = "  +4.58%  "
Set-TextValue $ws.Range("D39") "0.2146"
$ws.Range("E39").Value = "  -0.30%  "
Set-TextValue $ws.Range("D40") "0.6462"
$ws.Range("E40").Value = "  +2.17%  "
Set-TextValue $ws.Range("D41") "8.621"
$ws.Range("E41").Value = "  +0.96%  "
$ws.Range("E42").Value = "  +0.69%  "
$ws.Range("E43").Value = "  +0.69%  "
Set-TextValue $ws.Range("D44") "13.32"
$ws.Range("E44").Value = "  +2.95%  "
Set-TextValue $ws.Range("D45") "0.6042"
$ws.Range("E45").Value = "  +2.02%  "
Set-TextValue $ws.Range("D46") "2.165"
$ws.Range("E46").Value = "  +6.80%  "
Set-TextValue $ws.Range("D47") "3.620"
$ws.Range("E47").Value = "  -0.59%  "
Set-TextValue $ws.Range("D48") "122.17"
$ws.Range("E48").Value = "  -0.42%  "
Set-TextValue $ws.Range("D49") "1.207"
$ws.Range("E49").Value = "  -0.12%  "
Set-TextValue $ws.Range("D50") "1.137"
$ws.Range("E50").Value = "  +1.50%  "
Set-TextValue $ws.Range("D51") "77.89"
$ws.Range("E51").Value = "  +0.70%  "
